$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates: "Supports" -> "Accepts" in three cells, and year 2024 -> 2023 ---

$ws.Range("A2").Value = "Accepts tax on world top 1% to finance global poverty reduction`n(Additional 15% tax on income over [`$120k/year in PPP])"
$ws.Range("A3").Value = "Accepts tax on world top 3% to finance global poverty reduction`n(Additional 15% tax over [`$80k], 30% over [`$120k], 45% over [`$1M])"
$ws.Range("A8").Value = "Accepts reparations for colonization and slavery in`nthe form of funding education and technology transfers"
$ws.Range("A10").Value = '"My taxes ... global problems" (Global Nation, 2023)'

# --- Numeric updates for row 10 ---

$ws.Range("B10").Value = 0.55719177445442
$ws.Range("C10").Value = 0.590740460977192
$ws.Range("D10").Value = 0.43
$ws.Range("E10").Value = 0.65
$ws.Range("F10").Value = 0.76
$ws.Range("G10").Value = 0.58
$ws.Range("H10").Value = 0.6
$ws.Range("I10").Value = 0.52
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value = 0.76
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 0.44
